$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("About")
$ws2 = $wb.Worksheets.Item("PDiBCpDoC")

# -----------------------------------------------------------------
# "About" sheet: swap the citation from the MIT/RSC paper to BNEF
# -----------------------------------------------------------------
$ws.Range("B3").Value = "BNEF"
$ws.Range("B4").Value = 2019
$ws.Range("B5").Value = "https://about.bnef.com/blog/behind-scenes-take-lithium-ion-battery-prices/"
$ws.Range("B6").Value = "A Behind the Scenes Take on Lithium-ion Battery Prices"

# Drop the old "Abstract" pointer and the old learning-rate note - no
# longer applicable now that the source has changed.
$ws.Range("B7").ClearContents()
$ws.Range("A9").ClearContents()

# Add a new (still-empty) formatted cell further down the sheet - this
# is where the RPEF citation info will eventually live.
$ws.Range("D14").Interior.ColorIndex = -4142

# -----------------------------------------------------------------
# "PDiBCpDoC" sheet: hard-code the updated learning rate instead of
# averaging the two rates quoted in the old source
# -----------------------------------------------------------------
$ws2.Columns("B").ColumnWidth = 12.33
$ws2.Range("B2").Value = 0.18

# -----------------------------------------------------------------
# Restore on-screen selection / active sheet to match the new layout
# -----------------------------------------------------------------
$ws2.Activate()
$ws2.Range("H30").Select()

$ws.Activate()
$ws.Range("B11").Select()
